# The document's single table holds 100 "math fact" cells (20 rows x 5
# columns) read in row-major order. The commit replaces every cell's
# expression text with a new one (same position, same formatting) -
# update the whole grid in one pass rather than Find/Replace, since a
# couple of the old expressions (e.g. "72-33=") repeat and would
# otherwise collide.
$d = $word.ActiveDocument
$tbl = $d.Tables(1)

$newValues = @("34+55=", "12+62=", "16+18=", "30-19=", "50+19=", "73-32=", "56-29=", "21+13=", "49+12=", "43+37=", "4+65=", "17-3=", "89+10=", "11+66=", "27+28=", "64-25=", "20+24=", "71-3=", "36+59=", "42+5=", "83-55=", "66-25=", "4+50=", "2+67=", "3+83=", "65-55=", "39-10=", "25+12=", "38+10=", "15-2=", "98-79=", "6+87=", "0+80=", "94-83=", "38+38=", "80-31=", "66-33=", "36+20=", "56-36=", "22+59=", "20+3=", "29-0=", "25-3=", "71-61=", "8+88=", "21-7=", "21+78=", "64+24=", "85-70=", "89+1=", "49+2=", "62-30=", "16-15=", "80-54=", "25-14=", "76+2=", "7+88=", "92-19=", "18+0=", "54+5=", "85+2=", "52+45=", "54-39=", "22-19=", "51-38=", "80-11=", "21+60=", "37+35=", "27+43=", "1+92=", "64+21=", "5+41=", "55+16=", "10+77=", "7+33=", "87-21=", "25-11=", "20+44=", "83-19=", "17+78=", "90-61=", "78-32=", "52+42=", "19+52=", "67-57=", "55+31=", "9+75=", "22+6=", "12+53=", "90-41=", "58+21=", "5+4=", "95-71=", "60+34=", "23+13=", "35-29=", "16+29=", "2+75=", "55-9=", "17+75=")

$i = 0
foreach ($row in 1..$tbl.Rows.Count) {
    foreach ($col in 1..$tbl.Columns.Count) {
        $cell = $tbl.Cell($row, $col)
        $cell.Range.Text = $newValues[$i]
        $i = $i + 1
    }
}

Write-Host "Updated $i cells"
